# Auto-generated edit script applying cryptos.xlsx diff (Thu Apr 11 02:31:46 UTC 2024 GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.815.95"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "3.549.45"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").Value = "3.547.66"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +4.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.586"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "4.119.56"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "617.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.554.04"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.796.41"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.881"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -15.32%  "
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "574.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.101"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0469"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("D43").Value = "3.350.95"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.321"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "0.0₃0704"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "
